$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold font, border, centered alignment) from C1 to the new header cells D1:E1
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New header values for columns D and E (t+3, t+4)
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4

# Column B (rows 2-67) — recomputed predicted values
$Bvals = @(-0.3259533390584821,0.1978902336276767,0.1309485718638484,-0.08566793824230819,0.1276278515980622,-0.4614384876614183,-0.2709866285093435,-0.4209679496845873,0.3436887837841155,-0.2482459491696268,-0.1087508898535332,-0.01150366747922753,0.01155081633247283,-0.05835925389335367,0.285877318254868,0.56250432385674,0.04894586001739845,0.3804570174281816,0.2361495213532234,0.461768780998448,0.3414699679455382,-0.09238712393905801,4.534771994148857,0.524388496012612,0.3878503826166212,0.3341320112815402,1.096722682763318,5.863480215210426,1.00351129465251,-0.15159602179001,0.8266861627833791,0.9460996050162945,-0.6177245420987975,0.8228751327700211,0.7853821015207229,0.7609076149438667,0.7366364435688462,0.5855243524946664,0.7561083832988317,0.5622956694324531,0.717555762537827,0.7238315262564594,0.6827550122941594,0.6771136849576195,-1.259800347351256,-0.9744785939232784,-0.8670915112362984,-0.6387216172683954,-0.04741797340894434,-0.8597263583354644,-0.8597263583354644,-1.084026444314512,-0.1873877067204152,-0.9949076990006762,-0.8880174882996399,-0.9465298564717861,-1.150464319453582,-0.8597473676828168,-0.5169242166440466,0.3695974316258775,-1.227904722946455,-0.7638931319371256,-0.9076479785611221,-0.1299984289218459,-0.8131741492812394,-0.7940639888657253)
for ($i = 0; $i -lt $Bvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $Bvals[$i]
}

# Column C (rows 2-67) — recomputed predicted values
$Cvals = @(-0.3112759201248513,0.2068540516619886,0.1424338389046362,-0.07641264223041008,0.1432386406213986,-0.4522942319464852,-0.2592776833425081,-0.4111143807139855,0.3583482123923956,-0.2395354272017014,-0.0867907130281246,-0.01165389346983012,0.05155204817547648,-0.01410836412500775,0.3469334748427019,0.5975469357316218,0.03667612137647214,0.3958918984112226,0.2938092751582803,0.5186086745320461,0.3776735756792909,-0.06087923480789817,4.57384989143239,0.4894903269309043,0.3649107407635929,0.3022116221255643,1.066131646617195,5.447503492472463,0.9548851282734272,-0.2129154423199772,0.7871732809015993,0.9152751613758854,-0.6507460697090554,0.8138148514943915,0.7678300236883154,0.741357205273627,0.7174454321700235,0.581254185786493,0.7529239259480895,0.5542491414506645,0.6901709337605276,0.7089146264171076,0.6744044547749227,0.6583104991077348,-1.264221901809358,-0.9798357644125127,-0.8736975128145139,-0.642552574792616,-0.04890637255914625,-0.8653123345455569,-0.8653123345455569,-1.097711744551318,-0.1870535167737291,-0.9995363594736313,-0.8846002655561328,-0.9340088073217468,-1.12794034965257,-0.8480402148460283,-0.4971111734956855,0.3715007413611734,-1.208923604921344,-0.738311818883726,-0.9019909101174564,-0.1110487109336713,-0.7899089345158966,-0.7620704255356806)
for ($i = 0; $i -lt $Cvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $Cvals[$i]
}

# Column D (rows 2-67) — recomputed predicted values
$Dvals = @(-0.297031067946475,0.2167871348091792,0.1548873122604492,-0.06746381537253016,0.1583403482896266,-0.4425846283667368,-0.24675526487405,-0.398957653792802,0.3721432788652395,-0.2306230765166948,-0.06689733746444573,-0.01128321775599131,0.08797073075114029,0.02498140972080291,0.4017230581585344,0.6273381289153094,0.02605881697684193,0.4096285889324859,0.3446769913520271,0.5674358504612063,0.4089027974024245,-0.0340901860825687,4.530819290802055,0.4584631546021023,0.3443818550025777,0.2744108646979151,1.035108112914556,5.039878943785543,0.9101235324783095,-0.2640089275842031,0.7518820135477624,0.8883997727241983,-0.6806282358220184,0.8059150233966101,0.7529036636185144,0.7246334206242682,0.70071915508476,0.5780010875797771,0.7501121386842887,0.548426806697642,0.6660525226508269,0.696341770703753,0.667962365449804,0.6426328492272336,-1.268027398285052,-0.9844132606936714,-0.8786593157822313,-0.6451620197963193,-0.04993958249645333,-0.8695821122363085,-0.8695821122363085,-1.108985501860831,-0.1861402416707913,-1.003507091978611,-0.8821556501015857,-0.922904250613037,-1.10819690669765,-0.837058229996587,-0.4790572592892098,0.3739489402862353,-1.192233073768022,-0.7140604164949913,-0.8951714533026974,-0.09381236941148394,-0.7700727020603297,-0.7357373231785281)
for ($i = 0; $i -lt $Dvals.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $Dvals[$i]
}

# Column E (rows 2-67) — recomputed predicted values
$Evals = @(-0.2835478212216052,0.2267546620335315,0.167357836676539,-0.05931672983184233,0.1721059907394258,-0.4327701301594052,-0.2341453757227157,-0.3855131683035065,0.3844560053390241,-0.2220910555124325,-0.04954890393392026,-0.01090258543068098,0.1200001714382976,0.05835158744491686,0.4493533368504043,0.651411384363581,0.01645657032599904,0.4209024564469395,0.3880996029008814,0.6079756543920941,0.4346283633947823,-0.01233182337135161,4.421501485473843,0.4307215870663407,0.3232177437176343,0.2469851926634682,1.004102979255951,4.644718656198088,0.8684591437057199,-0.3094231159897771,0.7186382384955513,0.8638037386630153,-0.7080991982499313,0.798585027036251,0.7399609149213061,0.7100797884660245,0.685775935387294,0.5751782616211654,0.7471137381980206,0.5440186200575029,0.6446131761416967,0.6853766591219899,0.6627042221303827,0.6294800221062965,-1.271146924770325,-0.9882470557222144,-0.8821926283974424,-0.6468336844088383,-0.05092721814402677,-0.8726475335820014,-0.8726475335820014,-1.11800541456916,-0.1850279466862052,-1.006871462771747,-0.8806100719199673,-0.9132601035758863,-1.091082532953955,-0.8268423868653458,-0.4629535461656677,0.3762942870129932,-1.177809331338004,-0.6912917416230074,-0.8875848416232844,-0.07860579463933087,-0.7536157145179764,-0.714712208034043)
for ($i = 0; $i -lt $Evals.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $Evals[$i]
}

Write-Output "done"